# Auto-generated edit script: applies numeric market-data updates
# to the Siren_Profits-equivalent workbook (ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1532.5555
$ws.Range("I19").Value = 957
$ws.Range("J19").Value = 1820.3334
$ws.Range("K19").Value = 957
$ws.Range("L19").Value = 1820.3334
$ws.Range("M19").Value = -782
$ws.Range("N19").Value = -2170.3334

$ws.Range("H32").Value = 4700.6665
$ws.Range("J32").Value = 4700.6665
$ws.Range("L32").Value = 4700.6665
$ws.Range("N32").Value = -5352.6665

$ws.Range("H75").Value = 152999.25
$ws.Range("J75").Value = 152999.25
$ws.Range("L75").Value = 152999.25
$ws.Range("N75").Value = -154871.25

$ws.Range("H78").Value = 152999.25
$ws.Range("J78").Value = 152999.25
$ws.Range("L78").Value = 458997.75
$ws.Range("N78").Value = -468357.75

$ws.Range("H98").Value = 46036.645
$ws.Range("I98").Value = 51074.453
$ws.Range("J98").Value = 27564.666
$ws.Range("K98").Value = 51074.453
$ws.Range("L98").Value = 27564.666
$ws.Range("M98").Value = -49576.453
$ws.Range("N98").Value = -30560.666

$ws.Range("H112").Value = 73345.57000000001
$ws.Range("J112").Value = 78883.69500000001
$ws.Range("L112").Value = 236651.085
$ws.Range("N112").Value = -238867.085

$ws.Range("H116").Value = 1393383.1
$ws.Range("I116").Value = 11111111
$ws.Range("J116").Value = 5136.2856
$ws.Range("K116").Value = 11111111
$ws.Range("L116").Value = 5136.2856
$ws.Range("M116").Value = -11107669
$ws.Range("N116").Value = -12020.2856

$ws.Range("H122").Value = 46036.645
$ws.Range("I122").Value = 51074.453
$ws.Range("J122").Value = 27564.666
$ws.Range("K122").Value = 153223.359
$ws.Range("L122").Value = 82693.99800000001
$ws.Range("M122").Value = -150773.359
$ws.Range("N122").Value = -87593.99800000001

$ws.Range("H127").Value = 5720.8
$ws.Range("J127").Value = 7344.2856
$ws.Range("L127").Value = 22032.8568
$ws.Range("N127").Value = -31952.8568

$ws.Range("H138").Value = 3189.037
$ws.Range("J138").Value = 3844.4187
$ws.Range("L138").Value = 11533.2561
$ws.Range("N138").Value = -21813.2561

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 460
$ws.Range("I4").Value = 393
$ws.Range("K4").Value = 393
$ws.Range("M4").Value = -277

$ws.Range("H61").Value = 4141.8125
$ws.Range("I61").Value = 3404.4194
$ws.Range("K61").Value = 3404.4194
$ws.Range("M61").Value = -3192.4194

$ws.Range("H88").Value = 3112.25
$ws.Range("I88").Value = 2724.75
$ws.Range("K88").Value = 2724.75
$ws.Range("M88").Value = -2318.75

$ws.Range("H91").Value = 3112.25
$ws.Range("I91").Value = 2724.75
$ws.Range("K91").Value = 2724.75
$ws.Range("M91").Value = -1320.75

$ws.Range("H132").Value = 4918.1665
$ws.Range("I132").Value = 2261.125
$ws.Range("K132").Value = 6783.375
$ws.Range("M132").Value = -4253.375

$ws.Range("H136").Value = 4141.8125
$ws.Range("I136").Value = 3404.4194
$ws.Range("K136").Value = 10213.2582
$ws.Range("M136").Value = -7663.2582

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 80995.836
$ws.Range("J82").Value = 106743.75
$ws.Range("L82").Value = 106743.75
$ws.Range("N82").Value = -107509.75

$ws.Range("H85").Value = 80995.836
$ws.Range("J85").Value = 106743.75
$ws.Range("L85").Value = 106743.75
$ws.Range("N85").Value = -109395.75

$ws.Range("H92").Value = 40400.5
$ws.Range("J92").Value = 40400.5
$ws.Range("L92").Value = 40400.5
$ws.Range("N92").Value = -45392.5

$ws.Range("H134").Value = 4313.2
$ws.Range("I134").Value = 3171
$ws.Range("K134").Value = 9513
$ws.Range("M134").Value = -6978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H31").Value = 3047.1155
$ws.Range("I31").Value = 1740.3334
$ws.Range("J31").Value = 5987.375
$ws.Range("K31").Value = 1740.3334
$ws.Range("L31").Value = 5987.375
$ws.Range("M31").Value = -1445.3334
$ws.Range("N31").Value = -6577.375

$ws.Range("H34").Value = 3047.1155
$ws.Range("I34").Value = 1740.3334
$ws.Range("J34").Value = 5987.375
$ws.Range("K34").Value = 1740.3334
$ws.Range("L34").Value = 5987.375
$ws.Range("M34").Value = -1538.3334
$ws.Range("N34").Value = -6391.375

$ws.Range("H69").Value = 26077.6
$ws.Range("I69").Value = 20097
$ws.Range("K69").Value = 20097
$ws.Range("M69").Value = -19348

$ws.Range("H72").Value = 26077.6
$ws.Range("I72").Value = 20097
$ws.Range("K72").Value = 60291
$ws.Range("M72").Value = -56547

$ws.Range("H88").Value = 25083.416
$ws.Range("J88").Value = 23295.1
$ws.Range("L88").Value = 23295.1
$ws.Range("N88").Value = -24107.1

$ws.Range("H91").Value = 25083.416
$ws.Range("J91").Value = 23295.1
$ws.Range("L91").Value = 23295.1
$ws.Range("N91").Value = -26103.1

$ws.Range("H132").Value = 12566.648
$ws.Range("I132").Value = 1144.04
$ws.Range("K132").Value = 3432.12
$ws.Range("M132").Value = -902.1199999999999

$ws.Range("H134").Value = 3723.158
$ws.Range("I134").Value = 1809.8966
$ws.Range("K134").Value = 5429.6898
$ws.Range("M134").Value = -2894.6898

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 5462.5
$ws.Range("I41").Value = 925
$ws.Range("K41").Value = 925
$ws.Range("M41").Value = -570

$ws.Range("H75").Value = 89119
$ws.Range("I75").Value = 89119
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 89119
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -88245
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 89119
$ws.Range("I78").Value = 89119
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 267357
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -262989
$ws.Range("N78").ClearContents()

$ws.Range("H82").Value = 250000
$ws.Range("J82").Value = 250000
$ws.Range("L82").Value = 250000
$ws.Range("N82").Value = -250766

$ws.Range("H85").Value = 250000
$ws.Range("J85").Value = 250000
$ws.Range("L85").Value = 250000
$ws.Range("N85").Value = -252652

$ws.Range("H122").Value = 12099.723
$ws.Range("I122").Value = 10087.454
$ws.Range("J122").Value = 15261.857
$ws.Range("K122").Value = 30262.362
$ws.Range("L122").Value = 45785.571
$ws.Range("M122").Value = -27812.362
$ws.Range("N122").Value = -50685.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6386.8
$ws.Range("I22").Value = 5858.5
$ws.Range("J22").Value = 8500
$ws.Range("K22").Value = 5858.5
$ws.Range("L22").Value = 8500
$ws.Range("M22").Value = -5563.5
$ws.Range("N22").Value = -9090

$ws.Range("H27").Value = 6386.8
$ws.Range("I27").Value = 5858.5
$ws.Range("J27").Value = 8500
$ws.Range("K27").Value = 5858.5
$ws.Range("L27").Value = 8500
$ws.Range("M27").Value = -5751.5
$ws.Range("N27").Value = -8714

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H61").Value = 12568
$ws.Range("I61").Value = 13993.286
$ws.Range("K61").Value = 13993.286
$ws.Range("M61").Value = -13791.286

$ws.Range("H113").Value = 12568
$ws.Range("I113").Value = 13993.286
$ws.Range("K113").Value = 13993.286
$ws.Range("M113").Value = -11823.286

$ws.Range("H122").Value = 6517.5454
$ws.Range("I122").Value = 4673.5
$ws.Range("J122").Value = 7571.2856
$ws.Range("K122").Value = 14020.5
$ws.Range("L122").Value = 22713.8568
$ws.Range("M122").Value = -11570.5
$ws.Range("N122").Value = -27613.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 66666.664
$ws.Range("I64").Value = 66666.664
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 66666.664
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -66418.664
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 66666.664
$ws.Range("I67").Value = 66666.664
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 66666.664
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -65808.664
$ws.Range("N67").ClearContents()

$ws.Range("H122").Value = 23658.535
$ws.Range("I122").Value = 4610.6924
$ws.Range("J122").Value = 40166.668
$ws.Range("K122").Value = 13832.0772
$ws.Range("L122").Value = 120500.004
$ws.Range("M122").Value = -11382.0772
$ws.Range("N122").Value = -125400.004

